# Completed PRODUCT csv file for data import.
# Replace the Lorem-ipsum placeholder product names/prices/images with
# real product data, organized by category (3 products per category).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Category 1: RAM (rows 2-4) ---
$ws.Range("B2").Value = "2x4GB DDR3 Nature Valley"
$ws.Range("D2").Value = "ram.jpg"
$ws.Range("B3").Value = "2x8GB DDR3 KitKat"
$ws.Range("D3").Value = "ram.jpg"
$ws.Range("B4").Value = "2x8GB DDR4 Coffee Crisp"
$ws.Range("D4").Value = "ram.jpg"

# --- Category 2: Motherboards (rows 5-7) ---
$ws.Range("B5").Value = "Enterprise 4-Slot"
$ws.Range("D5").Value = "motherboard.jpg"
$ws.Range("B6").Value = "Nimitz 4-Slot"
$ws.Range("D6").Value = "motherboard.jpg"
$ws.Range("B7").Value = "Kitty Hawk 2-Slot"
$ws.Range("D7").Value = "motherboard.jpg"

# --- Category 3: GPUs (rows 8-10) ---
$ws.Range("B8").Value = "ATD Senna "
$ws.Range("D8").Value = "gpu.jpg"
$ws.Range("B9").Value = "GTZi Vettel"
$ws.Range("D9").Value = "gpu.jpg"
$ws.Range("B10").Value = "GTZ Clarkson"
$ws.Range("D10").Value = "gpu.jpg"

# --- Category 4: CPUs (rows 11-13) ---
$ws.Range("B11").Value = "Quad-core 4.0GHz Socrates"
$ws.Range("D11").Value = "cpu.jpg"
$ws.Range("B12").Value = "Quad-core 3.5GHz Plato"
$ws.Range("D12").Value = "cpu.jpg"
$ws.Range("B13").Value = "Quad-core 4.0GHz Nietzsche"
$ws.Range("D13").Value = "cpu.jpg"

# --- Category 5: Storage (rows 14-16) ---
$ws.Range("B14").Value = "Hummingbird 120GB SSD"
$ws.Range("D14").Value = "storage.jpg"
$ws.Range("B15").Value = "Pelican 2TB HDD"
$ws.Range("D15").Value = "storage.jpg"
$ws.Range("B16").Value = "Falcon 500GB SSD"
$ws.Range("C16").Value = 201.66
$ws.Range("D16").Value = "storage.jpg"

# --- Category 6: Power supplies (rows 17-19) ---
$ws.Range("B17").Value = "Franklin 650W"
$ws.Range("C17").Value = 83.2
$ws.Range("D17").Value = "psu.jpg"
$ws.Range("B18").Value = "Faraday 450W"
$ws.Range("D18").Value = "psu.jpg"
$ws.Range("B19").Value = "Tesla 650W"
$ws.Range("C19").Value = 100.94
$ws.Range("D19").Value = "psu.jpg"

# --- Category 7: Displays (rows 20-22) ---
$ws.Range("B20").Value = "Redeemer 24`" 1920x1080"
$ws.Range("D20").Value = "display.jpg"
$ws.Range("B21").Value = "Liberty 27`" 2560x1440"
$ws.Range("C21").Value = 403.88
$ws.Range("D21").Value = "display.jpg"
$ws.Range("B22").Value = "Rushmore 24`" 1920x1080"
$ws.Range("C22").Value = 343.1
$ws.Range("D22").Value = "display.jpg"

# --- Category 8: Keyboards (rows 23-25) ---
$ws.Range("B23").Value = "Changdao Cherry MX Red"
$ws.Range("C23").Value = 110.27
$ws.Range("D23").Value = "keyboard.jpg"
$ws.Range("B24").Value = "Claymore Cherry MX Brown"
$ws.Range("C24").Value = 150.58000000000001
$ws.Range("D24").Value = "keyboard.jpg"
$ws.Range("B25").Value = "Scimitar Cherry MX Red"
$ws.Range("C25").Value = 131.9
$ws.Range("D25").Value = "keyboard.jpg"

# --- Category 9: Mice (rows 26-28) ---
$ws.Range("B26").Value = "Rhino Optical Wired"
$ws.Range("C26").Value = 60.61
$ws.Range("D26").Value = "mouse.jpg"
$ws.Range("B27").Value = "Stag Optical Wired"
$ws.Range("D27").Value = "mouse.jpg"
$ws.Range("B28").Value = "Hercules Optical Wireless"
$ws.Range("C28").Value = 51.97
$ws.Range("D28").Value = "mouse.jpg"

# --- Category 10: Desktops (rows 29-31) ---
$ws.Range("B29").Value = "XenoWare Gaming PC"
$ws.Range("D29").Value = "desktop.jpg"
$ws.Range("B30").Value = "AZUZ Workstation PC"
$ws.Range("D30").Value = "desktop.jpg"
$ws.Range("B31").Value = "Pewlett Hackard Tower PC"
$ws.Range("D31").Value = "desktop.jpg"

# --- Category 11: Accessories (rows 32-34) ---
$ws.Range("B32").Value = "Blasterman Speakers"
$ws.Range("C32").Value = 32.61
$ws.Range("D32").Value = "accessory.jpg"
$ws.Range("B33").Value = "Big Brother Webcam"
$ws.Range("D33").Value = "accessory.jpg"
$ws.Range("B34").Value = "NSA Headphones"
$ws.Range("C34").Value = 52.29
$ws.Range("D34").Value = "accessory.jpg"

# Restore the cell selection highlighted by the author while editing (C16).
$ws.Range("C16").Select()
